$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("procedimientos")

# Add the new row of data (row 5)
$ws.Range("A5").Value = "proceso"
$ws.Range("B5").Value = "proceso.genera_rdc01"
$ws.Range("C5").Value = "Procedimiento para la generacion del RDC01"

# Column B width is no longer auto bestFit; set explicit custom width
# (Target stored width is 20.21875 characters; this runtime's ColumnWidth
# setter quantizes to 1/6-character steps, so 19.33 is the nearest input
# that lands on the closest achievable stored width, 20.1666..)
$ws.Columns.Item(2).ColumnWidth = 19.33

# Update the active selection to the new last cell, matching the diff
$ws.Range("B5").Select()
